$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.946.40"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.043.76"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'245.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'57.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").Value = "'0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "

$ws.Range("D12").Value = "'15.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "'0.877"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.04%  "

$ws.Range("D14").Value = "2.340.44"
$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").Value = "'5.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").Value = "2.041.65"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "'18.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.30%  "

$ws.Range("D18").Value = "36.923.65"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").Value = "'73.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").Value = "'235.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.73%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'169.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "

$ws.Range("D27").Value = "'2.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.86%  "

$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").Value = "'5.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.95%  "

$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("D31").Value = "'1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("D32").Value = "'4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.70%  "

$ws.Range("D33").Value = "'0.0612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'0.0865"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.89%  "

$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("D37").Value = "'2.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("E38").Value = "  -1.61%  "

$ws.Range("D39").Value = "'3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("E40").Value = "  +2.78%  "

$ws.Range("D41").Value = "'0.0982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.75%  "

$ws.Range("E42").Value = "  +1.11%  "

$ws.Range("E43").Value = "  +1.61%  "

$ws.Range("D44").Value = "'16.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("D45").Value = "'96.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").Value = "1.290.19"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("D47").Value = "'2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.35%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'3.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.98%  "

$ws.Range("D50").Value = "'6.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").Value = "2.224.37"
